$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 98.68130522039117
$ws.Range("G4").Value = 12.41175809502602
$ws.Range("H4").Value = 27738.68352655321
